# ============================================================================
# PlayerPerformance_4687.xlsx edit:
#  - insert a new "Player Info" sheet at the front
#  - rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" / "ODI Bowling"
#    and replace the howstat URL values with the bare match code
#  - drop a handful of stray empty INNING_NUMBER cells on "ODI Batting"
#  - append a new "ODI Batting Extra" sheet at the end
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlHAlignCenter
    $range.VerticalAlignment = -4160     # xlVAlignTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

function Convert-LinkColumnToCode($ws, [int]$col, [int]$firstRow, [int]$lastRow) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $text = $cell.Text
        if ($text -and $text.Contains("MatchCode=")) {
            $parts = $text -split "MatchCode="
            $code = $parts[1]
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

# ----------------------------------------------------------------------------
# 1) "ODI Batting" sheet (currently sheet #1)
# ----------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
Convert-LinkColumnToCode $batting 4 2 36

# stray empty INNING_NUMBER cells that should simply not exist
$staleInningCells = @(13, 27, 29, 30, 34)
foreach ($r in $staleInningCells) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# ----------------------------------------------------------------------------
# 2) "ODI Bowling" sheet (currently sheet #2)
# ----------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
Convert-LinkColumnToCode $bowling 2 2 31

# ----------------------------------------------------------------------------
# 3) New "Player Info" sheet, inserted before "ODI Batting"
# ----------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $playerInfoHeaders[$c - 1]
}
Set-HeaderStyle $playerInfo.Range("A1:D1")

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4687"
$playerInfo.Cells.Item(2, 2).Value = "Simi Singh"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Off Break"

# ----------------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, appended after "ODI Bowling"
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $battingExtra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
Set-HeaderStyle $battingExtra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraData = @(
    @("4391", 7,    "0", "0", "2.78%",  "NO"),
    @("4394", 7,    "3", "0", "14.35%", "NO"),
    @("4397", 7,    "0", "0", "4.93%",  "NO"),
    @("4426", $null, $null, $null, $null, "NO"),
    @("4427", 8,    "3", "0", "11.79%", "NO"),
    @("4442", 7,    "2", "2", "23.68%", "YES"),
    @("4444", 8,    "2", "0", "8.49%",  "NO"),
    @("4446", 8,    "1", "0", "5.02%",  "NO"),
    @("4448", 7,    "1", "0", "8.70%",  "NO"),
    @("4466", 7,    "4", "0", "23.20%", "NO"),
    @("4467", 7,    $null, $null, $null, "NO"),
    @("4468", 7,    "1", "0", "12.88%", "NO"),
    @("4474", $null, $null, $null, $null, "NO"),
    @("4475", 9,    $null, $null, $null, "NO"),
    @("4478", $null, $null, $null, $null, "NO"),
    @("4492", 6,    "0", "0", "1.32%",  "NO"),
    @("4494", $null, $null, $null, $null, "NO"),
    @("4496", 8,    $null, $null, $null, "NO"),
    @("4605", 9,    "4", "1", "10.00%", "NO"),
    @("4608", 8,    "0", "0", "7.41%",  "NO")
)

$row = 2
foreach ($rec in $extraData) {
    $codeCell = $battingExtra.Cells.Item($row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rec[0]

    if ($null -ne $rec[1]) {
        $battingExtra.Cells.Item($row, 2).Value = $rec[1]
    } else {
        $battingExtra.Cells.Item($row, 2).ClearContents()
    }

    for ($c = 3; $c -le 5; $c++) {
        $v = $rec[$c - 1]
        if ($null -ne $v) {
            $cell = $battingExtra.Cells.Item($row, $c)
            $cell.NumberFormat = "@"
            $cell.Value = $v
        } else {
            $battingExtra.Cells.Item($row, $c).ClearContents()
        }
    }

    $battingExtra.Cells.Item($row, 6).Value = $rec[5]
    $row++
}
